$d = $word.ActiveDocument

$pairs = @(
    @{old="803÷8="; new="915÷5="},
    @{old="106÷6="; new="110÷3="},
    @{old="110÷4="; new="486÷7="},
    @{old="804÷3="; new="482÷8="},
    @{old="313÷7="; new="681÷5="},
    @{old="708÷6="; new="104÷8="},
    @{old="795÷8="; new="599÷8="},
    @{old="474÷8="; new="913÷4="},
    @{old="226÷5="; new="482÷5="},
    @{old="511÷7="; new="400÷3="},
    @{old="163÷4="; new="659÷3="},
    @{old="183÷2="; new="468÷4="},
    @{old="691÷2="; new="252÷9="},
    @{old="267÷9="; new="120÷7="},
    @{old="579÷6="; new="778÷7="},
    @{old="744÷9="; new="140÷5="},
    @{old="792÷2="; new="678÷2="},
    @{old="569÷2="; new="920÷2="},
    @{old="167÷8="; new="216÷9="},
    @{old="123÷9="; new="480÷3="},
    @{old="680÷5="; new="119÷4="},
    @{old="526÷3="; new="992÷6="},
    @{old="712÷6="; new="447÷5="},
    @{old="332÷4="; new="529÷9="},
    @{old="888÷4="; new="391÷8="}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
